$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 122, shifting existing rows 122..157 down to 123..158
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row 122 with the new record's data
$ws.Range("A122").Value = 7
$ws.Range("B122").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C122").Value = "Ñuble"
$ws.Range("D122").Value = 45135
$ws.Range("D122").NumberFormat = $ws.Range("D123").NumberFormat
$ws.Range("E122").Value = 16
$ws.Range("F122").Value = 100112031
$ws.Range("G122").Value = "Poroto verde"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 80
$ws.Range("K122").Value = 20000
$ws.Range("L122").Value = 20000
$ws.Range("M122").Value = 20000
$ws.Range("N122").Value = "$/malla 25 kilos"
$ws.Range("O122").Value = "Perú"
$ws.Range("P122").Value = 800
$ws.Range("Q122").Value = 25
$ws.Range("R122").Value = "Hortaliza"
